# formating table parameter add asset
# Swap the "parameter" labels (column B) between rows 9-10 and rows 11-12,
# and update the active selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap parameter names so that:
#   B9  : "Remnant at Noon"      -> "Ship Draught#After"
#   B10 : "Sea Condition"        -> "Ship Draught#Fore"
#   B11 : "Ship Draught#After"   -> "Remnant at Noon"
#   B12 : "Ship Draught#Fore"    -> "Sea Condition"
$ws.Range("B9").Value  = "Ship Draught#After"
$ws.Range("B10").Value = "Ship Draught#Fore"
$ws.Range("B11").Value = "Remnant at Noon"
$ws.Range("B12").Value = "Sea Condition"

# Update the sheet's active selection to B9:B10 (active cell B9)
[void]$ws.Range("B9:B10").Select()
